# Apply "Last Fortnight" Rank / Days Won columns to the dashboard sheet,
# and refresh the current selection to match the new working cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns D ("Rank") and E ("Days Won") ---
# Clone the formatting of the existing header cell (C1, bold EB Garamond)
# onto D1:E1 so no new style/font entries are introduced, then set text.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Rank"
$ws.Range("E1").Value = "Days Won"

# --- Data rows (2-7): Last Fortnight Day Won / Rank values ---
# Clone the formatting of an existing data cell (C2, regular EB Garamond)
# onto D2:E7 first, then fill in the values.
$ws.Range("C2").Copy()
$ws.Range("D2:E7").PasteSpecial(-4122)

$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 0

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 4

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 8

$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0

# Clear the clipboard marquee left behind by Copy()
$excel.CutCopyMode = 0

# --- Update the active selection to match the post-edit state ---
$ws.Range("G9").Select()
